# Update the `average_county_temperature` column (AD) with refreshed
# NOAA-sourced county temperature values for the affected facility groups
# in NAICS 311513.
#
# Each contiguous block of rows below corresponds to a single facility_id
# (column H) whose average_county_temperature was recalculated from NOAA
# data; every row that belongs to a given facility receives the same
# updated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD12:AD16").Value   = 12.51681286549706
$ws.Range("AD17:AD21").Value   = 15.74228395061728
$ws.Range("AD22:AD61").Value   = 1.925925925925943
$ws.Range("AD62:AD71").Value   = 12.66820987654322
$ws.Range("AD77:AD81").Value   = -3.222222222222223
$ws.Range("AD114:AD118").Value = 1.925925925925943
$ws.Range("AD119:AD128").Value = 20.68981481481483
$ws.Range("AD129:AD138").Value = 14.96875
$ws.Range("AD139:AD163").Value = 1.925925925925943
$ws.Range("AD164:AD173").Value = -3.222222222222223
